$wb = $excel.ActiveWorkbook

# --- Step 1: Rename "Brian Miller" -> "Jennie Stewart" on the Users sheet ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B2").Value = "Jennie Stewart"

# --- Step 2: Add four new opportunity rows (Add Counterparty test data) ---
$wsAdd = $wb.Worksheets.Item("AddOpportunity")
$wsAdd.Range("A2:AB2").Copy($wsAdd.Range("A3:AB3"))
$wsAdd.Range("C3").Value = "Directs"
$wsAdd.Range("A2:AB2").Copy($wsAdd.Range("A4:AB4"))
$wsAdd.Range("C4").Value = "Primary Capital Advisory"
$wsAdd.Range("A2:AB2").Copy($wsAdd.Range("A5:AB5"))
$wsAdd.Range("C5").Value = "GP Stake Sale"
$wsAdd.Range("A2:AB2").Copy($wsAdd.Range("A6:AB6"))
$wsAdd.Range("C6").Value = "LP-Led Secondaries"

# --- Step 3: Add the new OppJobTypes lookup sheet, right after AddOpportunity ---
$wsJobTypes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAdd)
$wsJobTypes.Name = "OppJobTypes"
$wsJobTypes.Range("A1").Value = "JobTypes"
$wsJobTypes.Range("A1").Font.Bold = $true
$wsJobTypes.Range("A2").Value = "Activism Advisory"
$wsJobTypes.Range("A3").Value = "Directs"
$wsJobTypes.Range("A4").Value = "Primary Capital Advisory"
$wsJobTypes.Range("A5").Value = "GP Stake Sale"
$wsJobTypes.Range("A6").Value = "LP-Led Secondaries"

# --- Step 4: Reorder sheets - Users moves in front of AddOpportunity ---
# (re-fetch fresh references now that the sheet collection has changed)
$wsUsers = $wb.Worksheets.Item("Users")
$wsAdd = $wb.Worksheets.Item("AddOpportunity")
$wsUsers.Move($wsAdd)

# --- Step 5: restore per-sheet cell selections, leaving AddOpportunity active ---
$wb.Worksheets.Item("Users").Range("G6").Select()
$wb.Worksheets.Item("OppJobTypes").Range("E21").Select()
$wb.Worksheets.Item("AppName").Range("F22").Select()
$wb.Worksheets.Item("ModuleName").Range("G23").Select()
$wb.Worksheets.Item("AddOpportunity").Range("F19").Select()
